$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, shifting existing rows 149-276 down to 150-277
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new data
$ws.Range("A149").Value = 10
$ws.Range("B149").Value = "Vega Modelo de Temuco"
$ws.Range("C149").Value = "La Araucanía"
$ws.Range("D149").Value = 45072
$ws.Range("E149").Value = 9
$ws.Range("F149").Value = 100112012
$ws.Range("G149").Value = "Espinaca"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 30
$ws.Range("K149").Value = 9000
$ws.Range("L149").Value = 9000
$ws.Range("M149").Value = 9000
$ws.Range("N149").Value = "$/docena de atados"
$ws.Range("O149").Value = "Región de La Araucanía"
$ws.Range("P149").Value = 3000
$ws.Range("Q149").Value = 3
$ws.Range("R149").Value = "Hortaliza"
